$wb = $excel.ActiveWorkbook

$wsRef  = $wb.Worksheets.Item("reference traking")
$wsDist = $wb.Worksheets.Item("disturbance rejection")

# -----------------------------------------------------------------
# Update the "disturbance rejection" sheet values (columns B:E)
# for PID / IPD / DPI / PID_DIST results - "multiple equal poles 8"
# -----------------------------------------------------------------

# Row 2 - IAE
$wsDist.Range("B2").Value = 0.027006397956728607
$wsDist.Range("C2").Value = 0.026893981217712914
$wsDist.Range("D2").Value = 0.026903019585067682
$wsDist.Range("E2").Value = 0.004364121870257579

# Row 3 - Kp
$wsDist.Range("B3").Value = 30.019776248350702
$wsDist.Range("C3").Value = 30.018817954391618
$wsDist.Range("D3").Value = 30.011538791421234
$wsDist.Range("E3").Value = 99.956010938252788

# Row 4 - Ti
$wsDist.Range("B4").Value = 0.1074643643580597
$wsDist.Range("C4").Value = 0.10529127196098649
$wsDist.Range("D4").Value = 0.10551680732681155
$wsDist.Range("E4").Value = 0.10940963220711369

# Row 5 - Td
$wsDist.Range("B5").Value = 1.8299370463193918
$wsDist.Range("C5").Value = 1.84734396868327
$wsDist.Range("D5").Value = 1.8453924521164513
$wsDist.Range("E5").Value = 0.9965561804803903

# Row 6 - Ta
$wsDist.Range("E6").Value = 0.013734688769642182

# Row 9 - N
$wsDist.Range("B9").Value = 249.95797553755358
$wsDist.Range("C9").Value = 249.98667448699553
$wsDist.Range("D9").Value = 249.98836639623838
$wsDist.Range("E9").Value = 150.27676897481228

# Row 10 - alfa
$wsDist.Range("E10").Value = 17.676729830555718

# Row 12 - settlingtime
$wsDist.Range("B12").Value = 25.941207743364313
$wsDist.Range("C12").Value = 26.24061612171343
$wsDist.Range("D12").Value = 26.253627802783942
$wsDist.Range("E12").Value = 7.619290605319506

# Row 15 - peak
$wsDist.Range("B15").Value = 0.010734176645309927
$wsDist.Range("C15").Value = 0.010606216411124881
$wsDist.Range("D15").Value = 0.01062345372257986
$wsDist.Range("E15").Value = 0.0041618191671662684

# Row 16 - peaktime
$wsDist.Range("B16").Value = 0.66435766594882084
$wsDist.Range("C16").Value = 0.66138126033990485
$wsDist.Range("D16").Value = 0.66180835950768879
$wsDist.Range("E16").Value = 0.46492856340000921

# -----------------------------------------------------------------
# Widen column E on the "disturbance rejection" sheet
# -----------------------------------------------------------------
$wsDist.Columns.Item(5).ColumnWidth = 13.83

# -----------------------------------------------------------------
# Make "reference traking" the active / selected sheet (it was
# "disturbance rejection" before)
# -----------------------------------------------------------------
$wsRef.Activate()
$wsRef.Range("D22").Select()
